$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.708.39"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.589.32"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.40"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.24"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.815.82"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").Value = "1.577.38"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "27.692.17"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "220.07"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "0.0$([char]0x2083)0695"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("E22").Value = "  -4.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.60"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.17"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.86"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.15"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  -4.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").Value = "1.371.99"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.977"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.28"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("D47").Value = "1.726.16"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.81"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "0.0$([char]0x2086)0101"
$ws.Range("E49").Value = "  +10.78%  "
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  -1.05%  "
